# Generate Report for Handoff
# Updates the localization-status report: "b.md" moves from
# "Handed back: in sync with en-US" to "Ready for handoff" after a new
# handoff package was generated for it (zh-cn + de-de), and flags that the
# handback on file is stale relative to the newest source.

$wb = $excel.ActiveWorkbook

$readyForHandoff = 'Ready for handoff'
$genDate         = '2016-08-18 06:36:03'
$errorDetail     = 'The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/62dc30f41fb7807fe407ed06f0eed6cd4466f302/e2e/a.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/24b33f9c6d90ecca87518152bcddfe3ccc46ba8a/e2e/b.md.'

# --- Overview sheet: row 3 is the "b.md" row ---
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E3").Value = $readyForHandoff
$wsOverview.Range("F3").Value = $readyForHandoff
$wsOverview.Range("G3").Value = $genDate

# --- zh-cn sheet: row 3 is the "b.md" row ---
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C3").Value = $readyForHandoff
$wsZhCn.Range("F3").Value = "False"
$wsZhCn.Range("G3").Value = "b.63290e5768f688058c7b37413b0a5c26c308f864.zh-cn.xlf"
$wsZhCn.Range("H3").Value = "2016-08-18 06:35:56"
$wsZhCn.Range("P3").Value = $errorDetail
$wsZhCn.Columns.Item(16).ColumnWidth = 40

# --- de-de sheet: row 3 is the "b.md" row ---
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C3").Value = $readyForHandoff
$wsDeDe.Range("G3").Value = "b.63290e5768f688058c7b37413b0a5c26c308f864.de-de.xlf"
$wsDeDe.Range("H3").Value = $genDate
$wsDeDe.Range("P3").Value = $errorDetail
$wsDeDe.Columns.Item(16).ColumnWidth = 40
